# Updated cryptos list on Mon Nov 13 18:28:55 UTC 2023 with GitHub Actions
# Refreshes Price (col D) and Volume(1h) (col E) for each coin row; a few
# coins swapped rank position so their Coin/Link/Price/Volume cells moved
# to a different row. D-column values that look numeric are prefixed with
# a leading apostrophe so Excel stores them as text (matching the workbook's
# original text-formatted Price column, e.g. "36.748.63", "0.110").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.748.63'
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").Value = '2.088.45'
$ws.Range("E3").Value = '  +2.02%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '''245.21'
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").Value = '''0.652'
$ws.Range("E6").Value = '  -1.53%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '''54.21'
$ws.Range("E8").Value = '  -4.59%  '
$ws.Range("D9").Value = '''58.77'
$ws.Range("E9").Value = '  -2.09%  '
$ws.Range("D10").Value = '''0.366'
$ws.Range("E10").Value = '  -3.63%  '
$ws.Range("D11").Value = '''0.0762'
$ws.Range("E11").Value = '  -1.46%  '
$ws.Range("D12").Value = '''0.110'
$ws.Range("E12").Value = '  +1.12%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '''0.910'
$ws.Range("E13").Value = '  +5.45%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '''15.02'
$ws.Range("E14").Value = '  -4.59%  '
$ws.Range("D15").Value = '2.386.30'
$ws.Range("E15").Value = '  +1.74%  '
$ws.Range("D16").Value = '''5.50'
$ws.Range("E16").Value = '  -2.84%  '
$ws.Range("D17").Value = '2.058.78'
$ws.Range("E17").Value = '  +0.54%  '
$ws.Range("D18").Value = '36.707.21'
$ws.Range("E18").Value = '  -1.05%  '
$ws.Range("D19").Value = '''17.16'
$ws.Range("E19").Value = '  -3.63%  '
$ws.Range("D20").Value = '''72.59'
$ws.Range("E20").Value = '  -2.45%  '
$ws.Range("D21").Value = '0.0₃0881'
$ws.Range("E21").Value = '  -0.74%  '
$ws.Range("D22").Value = '''5.44'
$ws.Range("E22").Value = '  +1.74%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").Value = '''2.40'
$ws.Range("E25").Value = '  -2.64%  '
$ws.Range("D26").Value = '''9.81'
$ws.Range("E26").Value = '  +4.09%  '
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").Value = '''2.16'
$ws.Range("E27").Value = '  -0.59%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '''167.82'
$ws.Range("E28").Value = '  -0.63%  '
$ws.Range("D29").Value = '''20.64'
$ws.Range("E29").Value = '  +3.36%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '''5.33'
$ws.Range("E30").Value = '  +10.82%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '''0.123'
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("D32").Value = '''1.18'
$ws.Range("E32").Value = '  +5.06%  '
$ws.Range("D33").Value = '''4.73'
$ws.Range("E33").Value = '  +6.12%  '
$ws.Range("D34").Value = '''0.0609'
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("D35").Value = '''2.40'
$ws.Range("E35").Value = '  +7.05%  '
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("E37").Value = '  +3.83%  '
$ws.Range("D38").Value = '''0.0828'
$ws.Range("E38").Value = '  -6.48%  '
$ws.Range("D39").Value = '''1.27'
$ws.Range("E39").Value = '  -4.69%  '
$ws.Range("D40").Value = '''1.16'
$ws.Range("E40").Value = '  +2.29%  '
$ws.Range("D41").Value = '''0.0220'
$ws.Range("E41").Value = '  -0.66%  '
$ws.Range("D42").Value = '''0.0957'
$ws.Range("E42").Value = '  -2.40%  '
$ws.Range("D43").Value = '''4.82'
$ws.Range("E43").Value = '  -7.41%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '''96.27'
$ws.Range("E44").Value = '  +0.93%  '
$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D45").Value = '''2.85'
$ws.Range("E45").Value = '  -9.71%  '
$ws.Range("D46").Value = '''15.96'
$ws.Range("E46").Value = '  -6.79%  '
$ws.Range("D47").Value = '1.375.90'
$ws.Range("E47").Value = '  +8.80%  '
$ws.Range("D48").Value = '''7.28'
$ws.Range("E48").Value = '  +7.23%  '
$ws.Range("D49").Value = '''2.43'
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("D50").Value = '''2.91'
$ws.Range("E50").Value = '  +1.78%  '
$ws.Range("D51").Value = '2.269.55'
$ws.Range("E51").Value = '  +1.76%  '
